$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.297.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "'3.493.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'604.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("D6").Value = "'168.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.34%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").Value = "'3.493.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'6.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'0.575"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "'46.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "'4.056.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'8.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.76%  "
$ws.Range("D17").Value = "'608.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.71%  "
$ws.Range("D18").Value = "'3.498.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'69.356.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "'17.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "'10.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.73%  "
$ws.Range("D23").Value = "'0.873"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("D24").Value = "'15.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").Value = "'95.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("D26").Value = "'3.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'2.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").Value = "'9.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").Value = "'32.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'8.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("D32").Value = "'3.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").Value = "'6.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.17%  "
$ws.Range("D35").Value = "'552.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("D36").Value = "'10.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("D37").Value = "'3.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("D38").Value = "'56.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "

# Rows 39 and 40 swapped (FirstDigitalUSD now ranks above Hedera)
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.100"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.48%  "

$ws.Range("D41").Value = "'0.0445"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "'0.138"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'3.338.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "'0.323"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("D45").Value = "'32.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").Value = "'0.0₃0690"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "'2.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("D50").Value = "'134.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").Value = "'5.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.97%  "
